# Adds a new "BrentFuture" worksheet (after "BrentOilPrices") containing
# Brent oil future prices by month, mirroring the "ECB" / "BrentOilPrices"
# date+value layout, then makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet right after the last existing sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "BrentFuture"

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "FuturePrice"

# --- Monthly future-price data (date serial, price) ----------------------
$rows = @(
    @(43983, 25.57),
    @(44013, 29.21),
    @(44044, 31.51),
    @(44075, 33.33),
    @(44105, 34.75),
    @(44136, 35.71),
    @(44166, 36.99),
    @(44197, 37),
    @(44228, 37.55),
    @(44256, 38.07),
    @(44287, 38.45),
    @(44317, 38.81),
    @(44348, 39.14),
    @(44378, 39.4),
    @(44409, 39.69),
    @(44440, 39.98),
    @(44470, 40.27),
    @(44501, 40.56),
    @(44531, 40.86),
    @(44562, 41.12),
    @(44593, 41.4),
    @(44621, 41.7),
    @(44652, 42),
    @(44682, 42.3),
    @(44713, 42.59),
    @(44743, 42.81),
    @(44774, 43.03),
    @(44805, 43.24),
    @(44835, 43.45),
    @(44866, 43.66),
    @(44896, 43.86),
    @(44927, 44.09),
    @(44958, 44.33),
    @(44986, 44.57),
    @(45017, 44.82),
    @(45047, 45.06),
    @(45078, 45.3),
    @(45108, 45.52),
    @(45139, 45.72),
    @(45170, 45.9),
    @(45200, 46.08),
    @(45231, 46.26),
    @(45261, 46.41),
    @(45292, 46.63),
    @(45323, 46.84),
    @(45352, 47.04),
    @(45383, 47.22),
    @(45413, 47.4),
    @(45444, 47.58),
    @(45474, 47.76),
    @(45505, 47.94),
    @(45536, 48.12),
    @(45566, 48.3),
    @(45597, 48.48),
    @(45627, 48.65),
    @(45658, 48.83),
    @(45689, 49),
    @(45717, 49.17),
    @(45748, 49.33),
    @(45778, 49.49),
    @(45809, 49.65),
    @(45839, 49.81),
    @(45870, 49.97),
    @(45901, 50.13),
    @(45931, 50.29),
    @(45962, 50.45),
    @(45992, 50.61),
    @(46023, 50.76),
    @(46054, 50.91),
    @(46082, 51.06),
    @(46113, 51.21),
    @(46143, 51.38),
    @(46174, 51.55),
    @(46204, 51.72),
    @(46235, 51.89),
    @(46266, 52.06),
    @(46296, 52.23),
    @(46327, 52.4),
    @(46357, 52.57),
    @(46388, 52.63),
    @(46419, 52.69),
    @(46447, 52.75),
    @(46478, 52.81),
    @(46508, 52.87),
    @(46539, 52.93),
    @(46569, 52.99),
    @(46600, 53.05),
    @(46631, 53.11),
    @(46661, 53.15),
    @(46692, 53.2),
    @(46722, 53.25)
)

$dataArray = New-Object 'object[,]' $rows.Count,2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $dataArray[$i,0] = $rows[$i][0]
    $dataArray[$i,1] = $rows[$i][1]
}

$firstRow = 2
$lastRow = $firstRow + $rows.Count - 1
$ws.Range("A$firstRow`:B$lastRow").Value = $dataArray

# Date formatting for the populated date column plus a handful of blank
# trailing rows that were pre-formatted in the source sheet.
$ws.Range("A$firstRow`:A100").NumberFormat = "mmm-yy"

# --- Selection / active sheet state ---------------------------------------
$ws.Range("E89").Select()
